$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 117
$ws1.Range("F4").Value = 1625
$ws1.Range("F6").Value = 1096
$ws1.Range("F8").Value = 11522
$ws1.Range("F10").Value = 92
$ws1.Range("F11").Value = 451
$ws1.Range("F13").Value = 1092
$ws1.Range("F15").Value = 12385
$ws1.Range("F16").Value = 13094
$ws1.Range("F18").Value = 143
$ws1.Range("F21").Value = 228

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 118
$ws4.Range("F4").Value = 1625
$ws4.Range("F6").Value = 1096
$ws4.Range("F8").Value = 11522
$ws4.Range("F10").Value = 92
$ws4.Range("F11").Value = 451
$ws4.Range("F13").Value = 1092
$ws4.Range("F15").Value = 12385
$ws4.Range("F16").Value = 13095
$ws4.Range("F18").Value = 143
$ws4.Range("F21").Value = 228
